# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker table (B15:J18) lists one row per worker in arrears. A new
# worker record (LILIANA TERESA AYALA VILLALOBOS / CC 45553751, periodo
# mora 1707, valor mora 29509, salario basico 921696) is inserted ahead
# of the existing LAURA VELEZ URZOLA rows, and the second LAURA VELEZ
# URZOLA period (1911 / 48000 / 900000) moves down to the last row of
# the table. The first LAURA VELEZ URZOLA period (1910 / 41600 / 900000)
# stays put in row 17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 now holds the new worker (previously in row 18)
$ws.Range("C16").Value = "45553751"
$ws.Range("D16").Value = "LILIANA TERESA AYALA VILLALOBOS"
$ws.Range("E16").Value = "1707"
$ws.Range("F16").Value = 29509
$ws.Range("G16").Value = 921696

# Row 17 (LAURA VELEZ URZOLA, periodo 1910) is unchanged.

# Row 18 now holds LAURA VELEZ URZOLA's other period (previously in row 16)
$ws.Range("C18").Value = "1047401592"
$ws.Range("D18").Value = "LAURA VELEZ URZOLA"
$ws.Range("E18").Value = "1911"
$ws.Range("F18").Value = 48000
$ws.Range("G18").Value = 900000
